$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 53; this pushes the existing rows 53-103
# down to 55-105 (and therefore also appends two "new" rows at 104-105
# automatically, since the old 102-103 land there).
$ws.Rows("53:54").Insert()

# Populate the two freshly inserted rows (53 and 54) with their new data.
# Row 53: Cuatro cascos verde / Primera
$ws.Cells.Item(53, 1).Value2 = 12
$ws.Cells.Item(53, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(53, 3).Value = "Metropolitana"
$ws.Cells.Item(53, 4).Value2 = 45271
$ws.Cells.Item(53, 5).Value2 = 13
$ws.Cells.Item(53, 6).Value2 = 100112002
$ws.Cells.Item(53, 7).Value = "Pimiento"
$ws.Cells.Item(53, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value2 = 10
$ws.Cells.Item(53, 11).Value2 = 18000
$ws.Cells.Item(53, 12).Value2 = 18000
$ws.Cells.Item(53, 13).Value2 = 18000
$ws.Cells.Item(53, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(53, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(53, 16).Value2 = 1000
$ws.Cells.Item(53, 17).Value2 = 18
$ws.Cells.Item(53, 18).Value = "Hortaliza"

# Row 54: Cuatro cascos verde / Segunda
$ws.Cells.Item(54, 1).Value2 = 12
$ws.Cells.Item(54, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(54, 3).Value = "Metropolitana"
$ws.Cells.Item(54, 4).Value2 = 45271
$ws.Cells.Item(54, 5).Value2 = 13
$ws.Cells.Item(54, 6).Value2 = 100112002
$ws.Cells.Item(54, 7).Value = "Pimiento"
$ws.Cells.Item(54, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(54, 9).Value = "Segunda"
$ws.Cells.Item(54, 10).Value2 = 18
$ws.Cells.Item(54, 11).Value2 = 16000
$ws.Cells.Item(54, 12).Value2 = 16000
$ws.Cells.Item(54, 13).Value2 = 16000
$ws.Cells.Item(54, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(54, 16).Value2 = 889
$ws.Cells.Item(54, 17).Value2 = 18
$ws.Cells.Item(54, 18).Value = "Hortaliza"
